# Apply the "Fixed some bugs - Introduced category('Long') for tests" edit
# to the LanguageToolsTasks workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Row 15 - reword the Task Description cell (C15)
$ws.Range("C15").Value = "Count on number of rules, count on number of examples, total distance"

# 2) Row 22 - Progress % for "Deleting examples crashes the app" is now done (100%)
$ws.Range("E22").Value = 100

# 3) New row 24 - TransliterationEditor / Add Icon (Bug, Feature progress 0%)
$ws.Range("B24").Value = "TransliterationEditor"
$ws.Range("C24").Value = "Add Icon"
$ws.Range("D24").Value = "Feature"
$ws.Range("E24").Value = 0

# 4) New row 25 - TransliterationEditor / Adding rule or example with (') crashes the application (Bug, done)
$ws.Range("B25").Value = "TransliterationEditor"
$ws.Range("C25").Value = "Adding rule or example with (') crashes the application"
$ws.Range("D25").Value = "Bug"
$ws.Range("E25").Value = 100

# Match the formatting used by the other "Task Description" cells (wrap text)
$ws.Range("C24").WrapText = $true
$ws.Range("C25").WrapText = $true

# 5) Update the view: scrolled up one row, and the active selection moved to F17
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F17").Select()
